# Auto-generated Excel COM-interop edit script
# Applies the remaining-days decrement update (and the row-94 restock reset)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new value for column E (col 5), and optionally column F (col 6)
$updates = @(
    @{Row=2; E=11; F=$null}
    @{Row=3; E=11; F=$null}
    @{Row=4; E=11; F=$null}
    @{Row=5; E=5; F=$null}
    @{Row=6; E=11; F=$null}
    @{Row=7; E=5; F=$null}
    @{Row=8; E=11; F=$null}
    @{Row=9; E=5; F=$null}
    @{Row=10; E=4; F=$null}
    @{Row=11; E=11; F=$null}
    @{Row=12; E=5; F=$null}
    @{Row=13; E=11; F=$null}
    @{Row=14; E=11; F=$null}
    @{Row=15; E=11; F=$null}
    @{Row=16; E=9; F=$null}
    @{Row=17; E=5; F=$null}
    @{Row=18; E=8; F=$null}
    @{Row=19; E=8; F=$null}
    @{Row=20; E=8; F=$null}
    @{Row=21; E=8; F=$null}
    @{Row=22; E=5; F=$null}
    @{Row=23; E=5; F=$null}
    @{Row=24; E=5; F=$null}
    @{Row=25; E=5; F=$null}
    @{Row=26; E=5; F=$null}
    @{Row=27; E=5; F=$null}
    @{Row=28; E=8; F=$null}
    @{Row=29; E=8; F=$null}
    @{Row=30; E=8; F=$null}
    @{Row=31; E=8; F=$null}
    @{Row=32; E=8; F=$null}
    @{Row=33; E=8; F=$null}
    @{Row=34; E=8; F=$null}
    @{Row=35; E=8; F=$null}
    @{Row=37; E=8; F=$null}
    @{Row=38; E=8; F=$null}
    @{Row=39; E=8; F=$null}
    @{Row=40; E=4; F=$null}
    @{Row=41; E=4; F=$null}
    @{Row=42; E=8; F=$null}
    @{Row=43; E=5; F=$null}
    @{Row=44; E=4; F=$null}
    @{Row=45; E=5; F=$null}
    @{Row=46; E=4; F=$null}
    @{Row=47; E=8; F=$null}
    @{Row=48; E=4; F=$null}
    @{Row=49; E=5; F=$null}
    @{Row=50; E=3; F=$null}
    @{Row=51; E=3; F=$null}
    @{Row=52; E=3; F=$null}
    @{Row=53; E=3; F=$null}
    @{Row=54; E=3; F=$null}
    @{Row=55; E=3; F=$null}
    @{Row=56; E=3; F=$null}
    @{Row=57; E=3; F=$null}
    @{Row=58; E=7; F=$null}
    @{Row=59; E=7; F=$null}
    @{Row=60; E=7; F=$null}
    @{Row=61; E=5; F=$null}
    @{Row=62; E=7; F=$null}
    @{Row=63; E=7; F=$null}
    @{Row=64; E=7; F=$null}
    @{Row=65; E=8; F=$null}
    @{Row=66; E=8; F=$null}
    @{Row=67; E=8; F=$null}
    @{Row=68; E=8; F=$null}
    @{Row=69; E=8; F=$null}
    @{Row=70; E=9; F=$null}
    @{Row=71; E=9; F=$null}
    @{Row=72; E=9; F=$null}
    @{Row=73; E=9; F=$null}
    @{Row=74; E=9; F=$null}
    @{Row=75; E=9; F=$null}
    @{Row=76; E=9; F=$null}
    @{Row=77; E=2; F=$null}
    @{Row=78; E=2; F=$null}
    @{Row=79; E=2; F=$null}
    @{Row=80; E=2; F=$null}
    @{Row=81; E=2; F=$null}
    @{Row=82; E=2; F=$null}
    @{Row=83; E=2; F=$null}
    @{Row=84; E=2; F=$null}
    @{Row=85; E=2; F=$null}
    @{Row=86; E=2; F=$null}
    @{Row=87; E=4; F=$null}
    @{Row=88; E=4; F=$null}
    @{Row=89; E=4; F=$null}
    @{Row=90; E=4; F=$null}
    @{Row=91; E=5; F=$null}
    @{Row=92; E=4; F=$null}
    @{Row=93; E=2; F=$null}
    @{Row=94; E=7; F=20251120}
    @{Row=95; E=1; F=$null}
    @{Row=96; E=9; F=$null}
    @{Row=97; E=9; F=$null}
    @{Row=98; E=9; F=$null}
    @{Row=99; E=9; F=$null}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    if ($null -ne $u.F) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F
    }
}
